$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellRefs = @(
    'D2',
    'E2',
    'G2',
    'D3',
    'E3',
    'G3',
    'D4',
    'E4',
    'G4',
    'D5',
    'E5',
    'G5',
    'D6',
    'E6',
    'G6',
    'D7',
    'E7',
    'G7',
    'D8',
    'E8',
    'G8',
    'D9',
    'E9',
    'G9',
    'D10',
    'E10',
    'G10',
    'D11',
    'E11',
    'G11',
    'D12',
    'E12',
    'G12',
    'D13',
    'E13',
    'G13',
    'D14',
    'E14',
    'G14',
    'D15',
    'E15',
    'G15',
    'D16',
    'E16',
    'G16',
    'D17',
    'E17',
    'G17',
    'D18',
    'E18',
    'G18',
    'D19',
    'E19',
    'G19',
    'D20',
    'E20',
    'G20',
    'D21',
    'G21',
    'D22',
    'E22',
    'G22',
    'D23',
    'E23',
    'G23',
    'D24',
    'E24',
    'G24',
    'D25',
    'E25',
    'G25',
    'D26',
    'G26',
    'G27',
    'G28',
    'G29',
    'G30',
    'G31',
    'G32',
    'G33',
    'G34',
    'G35',
    'G36',
    'G37',
    'D38',
    'E38',
    'G38',
    'D39',
    'E39',
    'G39',
    'D40',
    'E40',
    'G40',
    'E41',
    'G41',
    'D42',
    'E42',
    'G42',
    'D43',
    'E43',
    'G43',
    'D44',
    'E44',
    'G44',
    'D45',
    'G45',
    'D46',
    'E46',
    'G46',
    'E47',
    'G47',
    'D48',
    'E48',
    'G48',
    'D49',
    'G49',
    'D50',
    'E50',
    'G50',
    'D51',
    'E51',
    'G51'
)

$newValues = @(
    '''327.31',
    '''-0.80%',
    '''20',
    '''44.30',
    '''0.73%',
    '''20',
    '''5.278',
    '''-4.42%',
    '''20',
    '''0.08349',
    '''3.01%',
    '''20',
    '''1.939',
    '''-6.34%',
    '''20',
    '''0.9731',
    '''-0.01%',
    '''20',
    '''2.499',
    '''-6.11%',
    '''20',
    '''0.1131',
    '''1.40%',
    '''20',
    '''0.1908',
    '''1.16%',
    '''20',
    '''0.09646',
    '''-3.43%',
    '''20',
    '''0.04613',
    '''-2.10%',
    '''20',
    '''0.1059',
    '''0.35%',
    '''20',
    '''0.001293',
    '''3.12%',
    '''20',
    '''0.006132',
    '''2.01%',
    '''20',
    '''3.403',
    '''1.80%',
    '''20',
    '''4.451',
    '''0.52%',
    '''20',
    '''0.3347',
    '''1.00%',
    '''20',
    '''8.709',
    '''-14.40%',
    '''20',
    '''0.1363',
    '''-1.93%',
    '''20',
    '''0.2580',
    '''20',
    '''0.04157',
    '''1.34%',
    '''20',
    '''0.001233',
    '''-5.75%',
    '''20',
    '''0.004406',
    '''0.29%',
    '''20',
    '''0.0001302',
    '''1.85%',
    '''20',
    '''0.0002985',
    '''20',
    '''20',
    '''20',
    '''20',
    '''20',
    '''20',
    '''20',
    '''20',
    '''20',
    '''20',
    '''20',
    '''20',
    '''0.02711',
    '''0.79%',
    '''20',
    '''0.05650',
    '''0.40%',
    '''20',
    '''0.007847',
    '''3.00%',
    '''20',
    '''-0.55%',
    '''20',
    '''0.007356',
    '''-3.26%',
    '''20',
    '''0.002120',
    '''8.34%',
    '''20',
    '''0.007859',
    '''-5.62%',
    '''20',
    '''0.3498',
    '''20',
    '''0.00006841',
    '''-3.50%',
    '''20',
    '''0.34%',
    '''20',
    '''0.003490',
    '''-2.03%',
    '''20',
    '''0.003538',
    '''20',
    '''0.00002105',
    '''0.34%',
    '''20',
    '''0.0002004',
    '''0.34%',
    '''20'
)

for ($i = 0; $i -lt $cellRefs.Count; $i++) {
    $cell = $ws.Range($cellRefs[$i])
    $cell.Value = $newValues[$i]
    $cell.Style = "Normal"
}
